{"js": "// Replace the division-problem text in each table cell according to the\n// fixed old-text -> new-text mapping below. We look the values up by their\n// current (pre-edit) text so the mapping is unambiguous even though a few\n// new values coincide with old values used elsewhere in the table\n// (e.g. \"10\u00f75=\" is both an original cell value and the replacement for a\n// different cell) \u2014 all originals are captured before any writes happen.\nconst replacements = {\n  \"90\u00f75=\": \"31\u00f76=\",\n  \"69\u00f72=\": \"10\u00f74=\",\n  \"41\u00f76=\": \"45\u00f76=\",\n  \"13\u00f73=\": \"64\u00f79=\",\n  \"30\u00f75=\": \"16\u00f73=\",\n  \"19\u00f73=\": \"43\u00f72=\",\n  \"89\u00f75=\": \"59\u00f74=\",\n  \"59\u00f75=\": \"32\u00f75=\",\n  \"97\u00f79=\": \"55\u00f73=\",\n  \"28\u00f79=\": \"53\u00f74=\",\n  \"31\u00f74=\": \"68\u00f74=\",\n  \"47\u00f76=\": \"86\u00f72=\",\n  \"10\u00f75=\": \"85\u00f76=\",\n  \"94\u00f79=\": \"99\u00f78=\",\n  \"19\u00f74=\": \"53\u00f78=\",\n  \"66\u00f74=\": \"60\u00f73=\",\n  \"80\u00f73=\": \"25\u00f76=\",\n  \"46\u00f74=\": \"10\u00f75=\",\n  \"21\u00f79=\": \"57\u00f75=\",\n  \"43\u00f78=\": \"85\u00f76=\",\n  \"39\u00f79=\": \"20\u00f78=\",\n  \"61\u00f75=\": \"13\u00f76=\",\n  \"16\u00f72=\": \"47\u00f75=\",\n  \"18\u00f78=\": \"24\u00f78=\",\n  \"89\u00f77=\": \"28\u00f72=\",\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Collect every cell together with its current value in one batch.\nconst cells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\n// Apply the replacement that matches each cell's original text. Using\n// `cell.value = ...` (Range.Text under the hood) keeps the existing\n// paragraph/run formatting (font, size, alignment) intact, unlike replacing\n// the whole cell body.\nfor (const cell of cells) {\n  const next = replacements[cell.value];\n  if (next !== undefined) {\n    cell.value = next;\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the division-problem text in each table cell according to the\n# fixed old-text -> new-text mapping below. Cells are matched by their\n# current (pre-edit) text, so the mapping stays correct even though a few\n# new values coincide with old values used elsewhere in the table\n# (e.g. \"10\u00f75=\" is both an original cell value and the replacement for a\n# different cell) - every cell is read before any of them are written.\n$replacements = @{\n  \"90\u00f75=\" = \"31\u00f76=\"\n  \"69\u00f72=\" = \"10\u00f74=\"\n  \"41\u00f76=\" = \"45\u00f76=\"\n  \"13\u00f73=\" = \"64\u00f79=\"\n  \"30\u00f75=\" = \"16\u00f73=\"\n  \"19\u00f73=\" = \"43\u00f72=\"\n  \"89\u00f75=\" = \"59\u00f74=\"\n  \"59\u00f75=\" = \"32\u00f75=\"\n  \"97\u00f79=\" = \"55\u00f73=\"\n  \"28\u00f79=\" = \"53\u00f74=\"\n  \"31\u00f74=\" = \"68\u00f74=\"\n  \"47\u00f76=\" = \"86\u00f72=\"\n  \"10\u00f75=\" = \"85\u00f76=\"\n  \"94\u00f79=\" = \"99\u00f78=\"\n  \"19\u00f74=\" = \"53\u00f78=\"\n  \"66\u00f74=\" = \"60\u00f73=\"\n  \"80\u00f73=\" = \"25\u00f76=\"\n  \"46\u00f74=\" = \"10\u00f75=\"\n  \"21\u00f79=\" = \"57\u00f75=\"\n  \"43\u00f78=\" = \"85\u00f76=\"\n  \"39\u00f79=\" = \"20\u00f78=\"\n  \"61\u00f75=\" = \"13\u00f76=\"\n  \"16\u00f72=\" = \"47\u00f75=\"\n  \"18\u00f78=\" = \"24\u00f78=\"\n  \"89\u00f77=\" = \"28\u00f72=\"\n}\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $table.Rows.Count; $r++) {\n  for ($c = 1; $c -le $table.Columns.Count; $c++) {\n    $cell = $table.Cell($r, $c)\n    $current = $cell.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($replacements.ContainsKey($current)) {\n      $cell.Range.Text = $replacements[$current]\n    }\n  }\n}\n"}
